# Auto-generated Excel COM-interop script
# Applies the cell value updates described by the commit "Add data for 2022-10-20"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("F2").Value = 74
$ws.Range("H2").Value = 94
$ws.Range("B3").Value = 70
$ws.Range("D3").Value = 115
$ws.Range("H3").Value = 124
$ws.Range("B6").Value = 324
$ws.Range("C6").Value = 413
$ws.Range("D6").Value = 349
$ws.Range("E6").Value = 378
$ws.Range("F6").Value = 444
$ws.Range("G6").Value = 398
$ws.Range("H6").Value = 390
$ws.Range("I6").Value = 437
$ws.Range("B7").Value = 439
$ws.Range("C7").Value = 552
$ws.Range("D7").Value = 549
$ws.Range("E7").Value = 569
$ws.Range("F7").Value = 634
$ws.Range("G7").Value = 599
$ws.Range("H7").Value = 622
$ws.Range("I7").Value = 733

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("B8").Value = 26
$ws.Range("D8").Value = 30
$ws.Range("E8").Value = 42
$ws.Range("G12").Value = 2
$ws.Range("H19").Value = 10
$ws.Range("D20").Value = 13
$ws.Range("I20").Value = 5
$ws.Range("B28").Value = 32
$ws.Range("C29").Value = 6
$ws.Range("B30").Value = 5
$ws.Range("D32").Value = 44
$ws.Range("E32").Value = 51
$ws.Range("F32").Value = 57
$ws.Range("F38").Value = 4
$ws.Range("C48").Value = 3
$ws.Range("E53").Value = 70
$ws.Range("H53").Value = 79
$ws.Range("B63").Value = 6
$ws.Range("F65").Value = 30
$ws.Range("H65").Value = 13
$ws.Range("F70").Value = 22
$ws.Range("H77").Value = 25
$ws.Range("H78").Value = 9
$ws.Range("H82").Value = 10
$ws.Range("G88").Value = 10
$ws.Range("B98").Value = 439
$ws.Range("C98").Value = 552
$ws.Range("D98").Value = 549
$ws.Range("E98").Value = 569
$ws.Range("F98").Value = 634
$ws.Range("G98").Value = 599
$ws.Range("H98").Value = 622
$ws.Range("I98").Value = 733

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("H2").Value = 5
$ws.Range("H7").Value = 25

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("B3").Value = 2
$ws.Range("B6").Value = 5

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("B5").Value = 17
$ws.Range("D5").Value = 20
$ws.Range("E5").Value = 32
$ws.Range("B6").Value = 26
$ws.Range("D6").Value = 30
$ws.Range("E6").Value = 42

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("D3").Value = 3
$ws.Range("I5").Value = 1
$ws.Range("D6").Value = 13
$ws.Range("I6").Value = 5

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("D6").Value = 31
$ws.Range("E6").Value = 40
$ws.Range("F6").Value = 48
$ws.Range("D7").Value = 44
$ws.Range("E7").Value = 51
$ws.Range("F7").Value = 57

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("G4").Value = 7
$ws.Range("G5").Value = 10

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("B6").Value = 28
$ws.Range("B7").Value = 32

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("H3").Value = 6
$ws.Range("H5").Value = 10

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("H4").Value = 5
$ws.Range("H5").Value = 9

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("H2").Value = 9
$ws.Range("E6").Value = 56
$ws.Range("E7").Value = 70
$ws.Range("H7").Value = 79

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("F2").Value = 2
$ws.Range("H5").Value = 6
$ws.Range("F6").Value = 30
$ws.Range("H6").Value = 13

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 6

$ws = $wb.Worksheets.Item("New City")
$ws.Range("B4").Value = 4
$ws.Range("B5").Value = 6

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("H3").Value = 2
$ws.Range("H6").Value = 10

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("C5").Value = 3
$ws.Range("C6").Value = 3

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("F2").Value = 2
$ws.Range("F5").Value = 22

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("E5").Value = 4
$ws.Range("E6").Value = 4

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 2
